# Update grouped results (Cosine Similarity / Euclidean Distance) across all
# "Study case" and "Experiment" sheets with newly generated values.

$wb = $excel.ActiveWorkbook

$updates = @{
    "Study case 1" = @{
        B2 = 0.7610360199449046
        C2 = 4.187896660234053
        B3 = 0.7488193612614225
        C3 = 4.274803858309955
    }
    "Study case 2" = @{
        B2 = 0.7512306351975814
        C2 = 4.850476179074347
        B3 = 0.73764237446588
        C3 = 4.951251305266382
    }
    "Study case 3" = @{
        B2 = 0.7610360199449046
        C2 = 4.187896660234053
        B3 = 0.7355651945656037
        C3 = 4.411217825065794
    }
    "Study case 4" = @{
        B2 = 0.7610360199449046
        C2 = 4.187896660234053
        B3 = 0.7355651945656037
        C3 = 4.411217825065794
    }
    "Experiment 1" = @{
        B2 = 0.6720536429424536
        C2 = 4.727530209210354
        B3 = 0.665094505510461
        C3 = 4.762326107438171
    }
    "Experiment 2" = @{
        B2 = 0.6720536429424536
        C2 = 4.727530209210354
        B3 = 0.6572713830079461
        C3 = 4.815911196452968
    }
    "Experiment 3" = @{
        B2 = 0.7082380710354689
        C2 = 5.442023715745888
        B3 = 0.6920513250057071
        C3 = 5.569715858163367
    }
    "Experiment 4" = @{
        B2 = 0.7082380710354689
        C2 = 5.442023715745888
        B3 = 0.7010033836736442
        C3 = 5.494419528748598
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cellValues = $updates[$sheetName]
    foreach ($addr in $cellValues.Keys) {
        $ws.Range($addr).Value = $cellValues[$addr]
    }
}
